$d = $word.ActiveDocument

# Locate the paragraph containing the red "DESCRIBIR..." comment and remove
# it entirely (including its paragraph mark).
$target = "DESCRIBIR QUE ENTORNO DE DESARROLLO SE ESTA UTILIZANDO, SI SE TIENE INSTALADO ALGUN PAQUETE/LIBRER" + [char]0x00CD + "A EN PARTICULAR PARA HACER QUE FUNCIONE EL SOFTWARE"

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $text = $r.Text
    if ($text -like "*$target*") {
        $r.Delete()
        break
    }
}

# The paragraph that used to follow the removed one (centered, empty) should
# now host a "_GoBack" bookmark (Word automatically recreates this bookmark
# at the last edit location; we add it explicitly here to match).
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -eq [char]0x000D -and $p.Alignment -eq 1) {
        $d.Bookmarks.Add("_GoBack", $r)
        break
    }
}
